$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header in G1 changes from "Date Created" to "Date Ordered"
$ws.Range("G1").Value = "Date Ordered"
